$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.91289999999996
$ws.Range("D4").Value = -7.653600000000004
$ws.Range("A6").Value = -22.64750000000002
$ws.Range("A7").Value = -20.32229999999997
$ws.Range("A8").Value = -22.60130000000002
$ws.Range("D8").Value = -8.599499999999997
$ws.Range("D9").Value = -7.441999999999994
$ws.Range("D12").Value = -6.027499999999995
$ws.Range("A16").Value = -22.17760000000002
$ws.Range("D17").Value = -8.34029999999999
$ws.Range("D18").Value = -9.285599999999997
$ws.Range("D19").Value = -8.537999999999993
$ws.Range("A20").Value = -20.44839999999999
$ws.Range("D20").Value = -7.029399999999995
$ws.Range("A21").Value = -20.24609999999998
$ws.Range("D26").Value = -7.288499999999999
$ws.Range("A28").Value = -21.85099999999999
$ws.Range("A29").Value = -21.20279999999998
$ws.Range("A30").Value = -21.8191
$ws.Range("D31").Value = -8.326399999999992
$ws.Range("A32").Value = -21.15009999999999
$ws.Range("D39").Value = -8.277299999999991
$ws.Range("A40").Value = -20.37890000000001
$ws.Range("D40").Value = -7.5805
$ws.Range("D41").Value = -8.374599999999988
$ws.Range("D42").Value = -8.367299999999991
$ws.Range("D43").Value = -7.526800000000005
$ws.Range("A46").Value = -21.83929999999999
$ws.Range("D47").Value = -7.590699999999998
$ws.Range("D48").Value = -7.700699999999997
$ws.Range("A51").Value = -21.61429999999999
$ws.Range("A52").Value = -22.2358
$ws.Range("D54").Value = -7.886999999999997
$ws.Range("A57").Value = -22.44210000000002
$ws.Range("A59").Value = -22.70310000000002
$ws.Range("A62").Value = -22.18550000000003
$ws.Range("D62").Value = -8.345599999999994
$ws.Range("D63").Value = -6.543799999999996
$ws.Range("D64").Value = -7.064199999999994
$ws.Range("A66").Value = -21.50020000000001
$ws.Range("A73").Value = -20.42269999999999
$ws.Range("A74").Value = -21.64229999999998
$ws.Range("D76").Value = -7.745500000000003
$ws.Range("A77").Value = -20.20799999999999
$ws.Range("D81").Value = -7.986500000000001
$ws.Range("D84").Value = -8.305899999999999
$ws.Range("D89").Value = -8.542900000000003
$ws.Range("A92").Value = -21.60640000000002
$ws.Range("D94").Value = -6.061499999999995
$ws.Range("A100").Value = -22.25970000000002
